# Update the "TrialsSetup" progress tracker values on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# REMASTER (CLOU) row (row 8): Progress goes from 0 to 12.5
$ws.Range("C8").Value = 12.5

# INNOVATE row (row 10): add a Progress value of 0
$ws.Range("C10").Value = 0

$wb.Save()
